$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns I..AQ (row 2 and row 3 share identical values)
$values = @{
    "I"  = 0.7203389830508474
    "J"  = 0.7203389830508474
    "K"  = 2.24
    "L"  = 0.3796610169491526
    "M"  = 1.8
    "N"  = 0.04054054054054054
    "O"  = 0.8035714285714285
    "P"  = 1.8
    "Q"  = 0.04054054054054054
    "R"  = 0.8035714285714285
    "U"  = 3.9
    "V"  = 0.08783783783783784
    "W"  = 0.06021505376344086
    "X"  = 0.02909367449159898
    "Y"  = 0.03112137927184188
    "Z"  = 0.188618925831202
    "AA" = 0.1358695652173913
    "AB" = 0.03694839371105331
    "AC" = 0.09892117150633799
    "AD" = 26.8
    "AF" = 26.8
    "AG" = 22.9
    "AH" = 0.3764044943820224
    "AI" = 0.4161490683229813
    "AJ" = 0.3402674591381873
    "AK" = 0.3785123966942149
    "AL" = 1.26
    "AM" = 1.26
    "AO" = 3.373015873015873
    "AQ" = 3.373015873015873
}

foreach ($row in 2, 3) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}
